# Apply cryptos.xlsx data refresh as described in the commit
# "Updated cryptos list on Fri May 31 00:42:10 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.309.38"
$ws.Range("E2").Value = "  +0.92%  "
$ws.Range("D3").Value = "3.746.89"
$ws.Range("E3").Value = "  -0.59%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.51%  "
$ws.Range("D7").Value = "3.745.99"
$ws.Range("E7").Value = "  -0.59%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  -1.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.159"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.45"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.445"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000257"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -8.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.63%  "
$ws.Range("D15").Value = "4.377.21"
$ws.Range("E15").Value = "  -0.64%  "
$ws.Range("D16").Value = "3.760.70"
$ws.Range("E16").Value = "  -0.22%  "
$ws.Range("D17").Value = "68.293.82"
$ws.Range("E17").Value = "  +0.93%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.80"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.34%  "
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "464.14"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.692"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000143"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.91%  "
$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.25%  "
$ws.Range("D30").Value = "3.894.74"
$ws.Range("E30").Value = "  -0.66%  "
$ws.Range("E31").Value = "  -5.39%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "29.87"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.87%  "
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.23"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.15"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.93%  "
$ws.Range("B35").Value = "Binance-PegBSC-USD"
$ws.Range("C35").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  --%  "
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.08"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.83%  "
$ws.Range("D37").Value = "3.701.22"
$ws.Range("E37").Value = "  -0.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.100"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.35"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -10.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.138"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.996"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.72%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.77"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.70%  "
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.301"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "43.27"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +9.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.48"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.90"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "45.95"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "146.75"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.67%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "386.54"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.31%  "
